$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the Description cell (C2) with the expanded DRA Landing Page test script,
# appending the additional verification steps separated by " || ".
$ws.Range("C2").Value = 'Verify that, accessing of the URL  takes the user to DRA application Landing page || Verify that DRA Landing page, displays application branding and logo || Verify that DRA Landing page, contains feature promotion and iconography in the marketing section || Verify that DRA Landing page, displays link to privacy statement and terms of use. || verify that DRA Landing page, displays the message and email id on the DRA landing page "Having trouble with sign-in? please contact DRA_support@thomsonreuters.com "'

# Resize row 2 to fit the now much longer wrapped text.
$ws.Rows.Item(2).RowHeight = 90

# Move the active selection to the cell that was edited.
$ws.Range("C2").Select()
